$wb = $excel.ActiveWorkbook

$design = $wb.Worksheets.Item("Design")
$design.Cells.Item(2, 1).Value = "sar"
$design.Cells.Item(2, 5).Value = "sar"
$design.Cells.Item(2, 9).Value = "asr"
$design.Cells.Item(2, 12).Value = "asr"
$design.Cells.Item(2, 17).Value = "sam"
$design.Cells.Item(2, 22).Value = 2
$design.Cells.Item(3, 1).Value = "presetet"
$design.Cells.Item(3, 4).Value = "pre"
$design.Cells.Item(3, 5).Value = "setet"
$design.Cells.Item(3, 8).Value = 8
$design.Cells.Item(3, 9).Value = "erpstete"
$design.Cells.Item(3, 11).Value = "erp"
$design.Cells.Item(3, 12).Value = "stete"
$design.Cells.Item(3, 16).Value = "sresetet"
$design.Cells.Item(3, 17).Value = "presotet"
$design.Cells.Item(3, 21).Value = 0
$design.Cells.Item(4, 1).Value = "seteter"
$design.Cells.Item(4, 5).Value = "setet"
$design.Cells.Item(4, 6).Value = "er"
$design.Cells.Item(4, 8).Value = 7
$design.Cells.Item(4, 9).Value = "etsetre"
$design.Cells.Item(4, 12).Value = "etset"
$design.Cells.Item(4, 13).Value = "re"
$design.Cells.Item(4, 17).Value = "teteter"
$design.Cells.Item(4, 18).Value = "setetor"
$design.Cells.Item(4, 22).Value = 0
$design.Cells.Item(4, 23).Value = 0
$design.Cells.Item(5, 1).Value = "desetetful"
$design.Cells.Item(5, 4).Value = "de"
$design.Cells.Item(5, 5).Value = "setet"
$design.Cells.Item(5, 6).Value = "ful"
$design.Cells.Item(5, 9).Value = "edtsetelfu"
$design.Cells.Item(5, 11).Value = "ed"
$design.Cells.Item(5, 12).Value = "tsete"
$design.Cells.Item(5, 13).Value = "lfu"
$design.Cells.Item(5, 16).Value = "dosetetful"
$design.Cells.Item(5, 17).Value = "degetetful"
$design.Cells.Item(5, 18).Value = "desetetfuy"
$design.Cells.Item(5, 21).Value = 1
$design.Cells.Item(5, 22).Value = 0
$design.Cells.Item(5, 23).Value = 2
$design.Cells.Item(6, 1).Value = "uninsut"
$design.Cells.Item(6, 3).Value = "un"
$design.Cells.Item(6, 4).Value = "in"
$design.Cells.Item(6, 5).Value = "sut"
$design.Cells.Item(6, 8).Value = 7
$design.Cells.Item(6, 9).Value = "nunitsu"
$design.Cells.Item(6, 10).Value = "nu"
$design.Cells.Item(6, 11).Value = "ni"
$design.Cells.Item(6, 12).Value = "tsu"
$design.Cells.Item(6, 15).Value = "ulinsut"
$design.Cells.Item(6, 16).Value = "unilsut"
$design.Cells.Item(6, 17).Value = "uninset"
$design.Cells.Item(6, 20).Value = 1
$design.Cells.Item(7, 1).Value = "owkiableful"
$design.Cells.Item(7, 5).Value = "owki"
$design.Cells.Item(7, 8).Value = 11
$design.Cells.Item(7, 9).Value = "ikowblaeulf"
$design.Cells.Item(7, 12).Value = "ikow"
$design.Cells.Item(7, 13).Value = "blae"
$design.Cells.Item(7, 14).Value = "ulf"
$design.Cells.Item(7, 17).Value = "iwkiableful"
$design.Cells.Item(7, 18).Value = "owkiebleful"
$design.Cells.Item(7, 19).Value = "owkiablefur"
$design.Cells.Item(7, 22).Value = 0
$design.Cells.Item(7, 23).Value = 0
$design.Cells.Item(7, 24).Value = 2
$design.Cells.Item(8, 1).Value = "deawtioner"
$design.Cells.Item(8, 5).Value = "awt"
$design.Cells.Item(8, 6).Value = "ion"
$design.Cells.Item(8, 7).Value = "er"
$design.Cells.Item(8, 8).Value = 10
$design.Cells.Item(8, 9).Value = "edwatoinre"
$design.Cells.Item(8, 12).Value = "wat"
$design.Cells.Item(8, 13).Value = "oin"
$design.Cells.Item(8, 14).Value = "re"
$design.Cells.Item(8, 16).Value = "seawtioner"
$design.Cells.Item(8, 17).Value = "deiwtioner"
$design.Cells.Item(8, 18).Value = "deawtiorer"
$design.Cells.Item(8, 19).Value = "deawtionew"
$design.Cells.Item(8, 21).Value = 0
$design.Cells.Item(8, 22).Value = 0
$design.Cells.Item(8, 23).Value = 2
$design.Cells.Item(8, 24).Value = 1
$design.Cells.Item(9, 1).Value = "predemetful"
$design.Cells.Item(9, 3).Value = "pre"
$design.Cells.Item(9, 4).Value = "de"
$design.Cells.Item(9, 5).Value = "met"
$design.Cells.Item(9, 6).Value = "ful"
$design.Cells.Item(9, 8).Value = 11
$design.Cells.Item(9, 9).Value = "rpeedtemlfu"
$design.Cells.Item(9, 10).Value = "rpe"
$design.Cells.Item(9, 11).Value = "ed"
$design.Cells.Item(9, 12).Value = "tem"
$design.Cells.Item(9, 13).Value = "lfu"
$design.Cells.Item(9, 15).Value = "tredemetful"
$design.Cells.Item(9, 16).Value = "predametful"
$design.Cells.Item(9, 17).Value = "predemotful"
$design.Cells.Item(9, 18).Value = "predemetfel"
$design.Cells.Item(9, 20).Value = 0
$design.Cells.Item(9, 21).Value = 1
$design.Cells.Item(9, 22).Value = 1
$design.Cells.Item(9, 23).Value = 1

$roots = $wb.Worksheets.Item("Roots")
$roots.Cells.Item(2, 1).Value = "sut"
$roots.Cells.Item(3, 1).Value = "upi"
$roots.Cells.Item(4, 1).Value = "boha"
$roots.Cells.Item(5, 1).Value = "eget"
$roots.Cells.Item(6, 1).Value = "setet"
$roots.Cells.Item(7, 1).Value = "owki"
$roots.Cells.Item(8, 1).Value = "awt"
$roots.Cells.Item(9, 1).Value = "rers"
$roots.Cells.Item(10, 1).Value = "met"
$roots.Cells.Item(11, 1).Value = "sar"
